$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.763.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.540.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.09%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.540.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("E10").Value = "  -0.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.996.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.796.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.543.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("E25").Value = "  +7.72%  "

$ws.Range("E26").Value = "  -3.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0808"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "178.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "408.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.50%  "

$ws.Range("E40").Value = "  +3.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("E49").Value = "  +4.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("E51").Value = "  +0.07%  "
